$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style (bold, border, centered) of the
# other header cells by copying the format from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the new "Save" column data values (0) for rows 2 and 3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
